# "alteracao ciclo de vida AC6"
# The deck's date placeholder (an auto-updating datetimeFigureOut field,
# cached as "17/05/2019") is refreshed to "18/05/2019" on the slide
# master and on every slide layout that carries that placeholder.

$p = $ppt.ActivePresentation
$oldDate = "17/05/2019"
$newDate = "18/05/2019"

function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster

# Slide master's own date placeholder.
Update-DatePlaceholder $master.Shapes

# Every custom (slide) layout under the master has its own copy of the
# date placeholder too.
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholder $layout.Shapes
}
